$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new (empty) column before E. This shifts the existing
#    link column E -> F, carrying over its values/styles intact, and
#    moves the sheet's used range out to F.
# ------------------------------------------------------------------
$ws.Columns("E:E").Insert()

# ------------------------------------------------------------------
# Stash two donor cells (off in scratch column H) holding the exact
# "hyperlink" look-and-feel (cellXf index 2 / index 3) the F column
# already uses, *before* touching the Hyperlinks collection below --
# Hyperlinks.Add() stamps its own fresh style on whatever cell it
# touches, so the originals can't be trusted as donors afterward.
# ------------------------------------------------------------------
$ws.Range("F2").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Row 33 gains a "stack" tag in column B.
# ------------------------------------------------------------------
$ws.Range("B33").Value2 = "stack"

# ------------------------------------------------------------------
# 3. Seven new exercise rows (35-41).
# ------------------------------------------------------------------
$ws.Range("A35").Value2 = "implement trie prefix tree"
$ws.Range("B35").Value2 = "string"
$ws.Range("C35").Value2 = "dict of dicts"
$ws.Range("D35").Value2 = "tree"
$ws.Range("F35").Value2 = "https://leetcode.com/problems/implement-trie-prefix-tree/"

$ws.Range("A36").Value2 = "coin change"
$ws.Range("B36").Value2 = "dynamic programming"
$ws.Range("F36").Value2 = "https://leetcode.com/problems/coin-change"

$ws.Range("A37").Value2 = "product of array except self"
$ws.Range("B37").Value2 = "array"
$ws.Range("C37").Value2 = "two pointers"
$ws.Range("D37").Value2 = "cum prod"
$ws.Range("F37").Value2 = "https://leetcode.com/problems/product-of-array-except-self/"

$ws.Range("A38").Value2 = "min stack"
$ws.Range("B38").Value2 = "two stacks"
$ws.Range("F38").Value2 = "https://leetcode.com/problems/min-stack/"

$ws.Range("A39").Value2 = "validate binary search tree"
$ws.Range("B39").Value2 = "recursive"
$ws.Range("C39").Value2 = "dfs"
$ws.Range("D39").Value2 = "tree"
$ws.Range("F39").Value2 = "https://leetcode.com/problems/validate-binary-search-tree/"

$ws.Range("A40").Value2 = "number of islands"
$ws.Range("B40").Value2 = "recursive"
$ws.Range("C40").Value2 = "dfs"
$ws.Range("D40").Value2 = "matrix"
$ws.Range("E40").Value2 = 24
$ws.Range("F40").Value2 = "https://leetcode.com/problems/number-of-islands/"

$ws.Range("A41").Value2 = "rotting oranges"
$ws.Range("B41").Value2 = "stack"
$ws.Range("C41").Value2 = "bfs"
$ws.Range("D41").Value2 = "matrix"
$ws.Range("E41").Value2 = 26
$ws.Range("F41").Value2 = "https://leetcode.com/problems/rotting-oranges/"

# ------------------------------------------------------------------
# 4. Re-create the hyperlinks. The column insert above does not carry
#    the Hyperlinks collection's anchors from E to F, so drop them all
#    and re-add in the original order (keeps rId1..rId8 identical),
#    then append the two new links (rId9, rId10).
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://leetcode.com/problems/two-sum/")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://leetcode.com/problems/valid-parentheses/")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://leetcode.com/problems/valid-palindrome/")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://leetcode.com/problems/merge-two-sorted-lists/")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://leetcode.com/problems/valid-anagram/")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://leetcode.com/problems/binary-search/")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://leetcode.com/problems/flood-fill/")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://leetcode.com/problems/lowest-common-ancestor-of-a-binary-tree/")
$ws.Hyperlinks.Add($ws.Range("F36"), "https://leetcode.com/problems/coin-change")
$ws.Hyperlinks.Add($ws.Range("F40"), "https://leetcode.com/problems/number-of-islands/")

# ------------------------------------------------------------------
# 5. Put the original styling back on every cell Hyperlinks.Add()
#    just re-stamped, plus style the three brand-new non-hyperlink
#    "accent" cells (F35, F37, F39) the same way F1/F6/F31-34 look.
# ------------------------------------------------------------------
$ws.Range("H1").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F36").PasteSpecial(-4122)
$ws.Range("F40").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F10").PasteSpecial(-4122)

$ws.Range("F1").Copy()
$ws.Range("F35").PasteSpecial(-4122)
$ws.Range("F37").PasteSpecial(-4122)
$ws.Range("F39").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Drop the scratch donor cells so they don't linger in the saved sheet.
$ws.Range("H1:H2").Clear()

# ------------------------------------------------------------------
# 6. Move the selection the same way the source workbook does (now
#    pointing just past the newly-added rows).
# ------------------------------------------------------------------
$ws.Range("D42").Select()
